# Rename the two sheets to drop the "gi|...|ref|...|" wrapper, keeping just
# the NC_ accession.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "NC_017250.1"
$ws2.Name = "NC_017251.1"

# Restore/keep the selection on sheet 1 (J2), without it being the tab shown
# when the workbook is opened.
$ws1.Range("J2").Select()

# Sheet 2 becomes the active (visible) tab, with its selection moved to A37.
$ws2.Activate()
$ws2.Range("A37").Select()
